$wb = $excel.ActiveWorkbook

# The dated sheets hold per-day attendance/vitals readings; the sheet with the
# face-box timing fix is "2021-01-06" (it is not the active tab), so select it
# explicitly by name rather than relying on $wb.ActiveSheet.
$ws = $wb.Worksheets.Item("2021-01-06")

# Newly captured readings (face-detection box fix let these register) to
# append below the existing 4 data rows (rows 2-5).
$newRows = @(
    @{ Row = 6; Time = "14:18:15"; SpO2 = 96.6712521745384;  HR = 56.68526035361902 },
    @{ Row = 7; Time = "14:20:23"; SpO2 = 97.15096348632157; HR = 77.93673309236452 },
    @{ Row = 8; Time = "14:20:42"; SpO2 = 94.97829344384658; HR = 80.72980249297832 }
)

foreach ($r in $newRows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value = 1
    $ws.Cells.Item($i, 2).Value = "sachin"
    $ws.Cells.Item($i, 3).Value = "301/Sanskruti-1,Andheri, Mumbai"
    $ws.Cells.Item($i, 4).Value = "Software Engineer"
    $ws.Cells.Item($i, 5).Value = $r.Time
    $ws.Cells.Item($i, 6).Value = $r.SpO2
    $ws.Cells.Item($i, 7).Value = $r.HR
    $ws.Cells.Item($i, 8).Value = "NA"
    $ws.Cells.Item($i, 9).Value = "NA"

    # Mirror the bold, thin-bordered, centered look used for the "Sr. No"
    # column (A) in the existing data rows.
    $cell = $ws.Cells.Item($i, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}
